# 20150113 +++++++ cs-厂商 end ++ end
# Append the "我的订单" (my-indent) rows to the page-mapping sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A entries first (so the shared-string table fills in the same
# order as the authored workbook: A20, A21, then B20, B21).
$ws.Range("A20").Value = "我的订单3.psd,我的订单3 - 下拉.psd,我的订单-待发货.psd,"
$ws.Range("A21").Value = "我的订单-待发货-订单详情 - 申请退款.psd,我的订单-待发货-订单详情 - 退款成功.psd,我的订单-待发货-订单详情.psd"
$ws.Range("B20").Value = "my-indent.html"
$ws.Range("B21").Value = "my-indent-dfh.html"

# The PSD-file column wraps onto several lines, so give it the wrap-text
# style (this becomes the new 4th cellXfs entry) and size the rows to fit.
$ws.Range("A20:A21").WrapText = $true
$ws.Rows.Item(20).RowHeight = 27
$ws.Rows.Item(21).RowHeight = 40.5

# Match the author's final selection.
$ws.Range("F19").Select() | Out-Null
